$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 3 - this shifts the existing data
# (old rows 3,4,5) down to rows 4,5,6, matching the target layout.
[void]$ws.Rows.Item(3).Insert()

# New "Account" value for the (now) row 4 transaction.
$ws.Cells.Item(4, 3).Value = 70

# New trailing row with a date-formatted (but otherwise blank) cell in A7,
# mirroring the style used for the date column elsewhere on the sheet.
$ws.Cells.Item(7, 1).NumberFormat = "YYYY\-MM\-DD"

# Move/restore the active selection.
[void]$ws.Range("E9").Select()
